$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "По паспарту" -> "По паспорту"
$ws.Range("B2").Value = "По паспорту"

# New "Options" column (F) — schema change for printed ambulatory card form
$ws.Range("F1").Value = "Options"
$ws.Range("F2").Value = "|forAmbCard|"
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
